$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 81, shifting existing rows 81-126 down to 82-127.
$ws.Rows.Item(81).Insert()

# Fill in the newly inserted row 81 with a new weekly price record.
$ws.Cells.Item(81, 1).Value = 9
$ws.Cells.Item(81, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 44523
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100101
$ws.Cells.Item(81, 8).Value = "Berries"
$ws.Cells.Item(81, 9).Value = 100101001
$ws.Cells.Item(81, 10).Value = "Arándano (blue)"
$ws.Cells.Item(81, 11).Value = "Sin especificar"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 500
$ws.Cells.Item(81, 14).Value = 5000
$ws.Cells.Item(81, 15).Value = 5000
$ws.Cells.Item(81, 16).Value = 5000
$ws.Cells.Item(81, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(81, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(81, 19).Value = 2500
$ws.Cells.Item(81, 20).Value = 2

# Make sure the date cell keeps the same date number-format as the rest of column D.
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(82, 4).NumberFormat
